# DTFS2-8124: accrual currency validation changed to ISO4217
# Adds "Accrual currency" / "accrual exchange rate" columns (L, M) with
# sample invalid-currency data rows used by the e2e test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("L1").Value = "Accrual currency"
$ws.Range("M1").Value = "accrual exchange rate"

# Row 2 sample data (invalid accrual currency code)
$ws.Range("L2").Value = "INRA"
$ws.Range("M2").Value = 1.223

# Row 3 sample data (invalid accrual currency code)
$ws.Range("L3").Value = "A"
$ws.Range("M3").Value = 2.33

# Match the author's saved selection state
$ws.Range("L1:M3").Select()
